$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New admin accounts added under the existing Username/Password table.
$ws.Range("A3").Value = "Admin2"
$ws.Range("B3").Value = "admin1234"
$ws.Range("A4").Value = "Admin3"
$ws.Range("B4").Value = "admin12345"

# Let the two data columns fit their (now longer) contents.
$ws.Columns("A:B").AutoFit() | Out-Null

# Leave the cursor where the author's last save left it.
$ws.Range("C17").Select() | Out-Null
